# Fix Training Data Issue (#48)
# The "Date" column (BF) was populated with the source workbook's file-name
# stamp ("2-26-2011-12") instead of an actual date. Because of how the NBA
# stats were captured, the data was also off by a day; the correct value is
# the ISO date 2012-02-26.
#
# Rows 2-31 all hold the same bad literal, so walk them and fix each cell.
# NOTE: writing the ISO-looking string straight into Value/Value2 causes
# Excel to auto-coerce it into a date serial number (and stamp a date
# NumberFormat on the cell, which would also bump its style index). To keep
# the cell a genuine text value with its original (default) style, we enter
# it as a formula that evaluates to the literal string, then immediately
# convert that formula to a static value via copy / paste-special-values
# (xlPasteValues = -4163).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDate = "2-26-2011-12"
$newDate = "2012-02-26"

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    if ($cell.Value2 -eq $oldDate) {
        $cell.Formula = '="' + $newDate + '"'
        $cell.Copy()
        $cell.PasteSpecial(-4163)
    }
}
